# Improve analysis of 103124 and 110724
# Zero out Num_Inclusions (D), Inclusion_Area (E) and Overlap_Area (G) values
# for rows that should no longer report inclusions, for the rows affected
# by the re-analysis (rows 10, 67, 97, 99 on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: D10 -> 0, E10 -> 0, G10 -> 0 (F10 unchanged)
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 0

# Row 67: D67 -> 0, E67 -> 0 (F67, G67 unchanged)
$ws.Range("D67").Value = 0
$ws.Range("E67").Value = 0

# Row 97: D97 -> 0, E97 -> 0 (F97, G97 unchanged)
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0

# Row 99: D99 -> 0, E99 -> 0, G99 -> 0 (F99 unchanged)
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("G99").Value = 0
